$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhCn = $wb.Worksheets.Item("zh-cn")
$deDe = $wb.Worksheets.Item("de-de")

# Row 3 on every sheet is the "b5d4e24c-0b11-489b-a3a2-aba74e5ebd03" entry.
# The handback transform failed for this file on both locales, so its status
# flips from "Ready for handoff" to "Handback transform failed" everywhere it
# is shown (Overview summary columns + each locale's Status column), and the
# locale sheets' Error Detail column (column P) is populated with the reason.

$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

$zhCn.Range("C3").Value = "Handback transform failed"
$zhCn.Range("P3").Value = "Handback file name: lwy3iv0b.mij is different with handoff file name: b5d4e24c-0b11-489b-a3a2-aba74e5ebd03.7d675b57ec573fc9664b218021f4abdaef3188f9.zh-cn."
$zhCn.Columns.Item(16).ColumnWidth = 39.16666666666667

$deDe.Range("C3").Value = "Handback transform failed"
$deDe.Range("P3").Value = "Handback file name: lwy3iv0b.mij is different with handoff file name: b5d4e24c-0b11-489b-a3a2-aba74e5ebd03.7d675b57ec573fc9664b218021f4abdaef3188f9.de-de."
$deDe.Columns.Item(16).ColumnWidth = 39.16666666666667
